# Working Hours.xlsx update: "fixed homepage navbar and added custom scroll bar and table css"
#
# Husein (column C / H / J) logged work for several more days:
#   Row 30 (2019-11-05): hours corrected 3 -> 4
#   Row 31 (2019-11-06): 5 hrs, "shake effect and table", 3 tasks
#   Row 32 (2019-11-07): 6.5 hrs, "homepage navbar", 1 task
#   Row 34 (2019-11-09): 7 hrs, "fixed homepage navbar and custom scroll bar and table css", 1 task
#   Row 35 (2019-11-10): 5 hrs, "Fixed navbar in faculty login, removed navbar & set content to
#                          center in main login page", 3 tasks
# Column H (task text) is widened to fit the longer descriptions, and the view
# scrolls down to the newly active area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: correct Husein's hours for that day ---
$ws.Range("C30").Value = 4

# Task text is entered in the order it was typed in the original editing
# session (row 35's note first) so new shared-string entries land in the
# same order as the authoritative workbook.
$ws.Range("H35").Value = "Fixed navbar in faculty login, removed navbar & set content to center in main login page"
$ws.Range("H31").Value = "shake effect and table"
$ws.Range("H32").Value = "homepage navbar"
$ws.Range("H34").Value = "fixed homepage navbar and custom scroll bar and table css"

# --- Row 31 ---
$ws.Range("C31").Value = 5
$ws.Range("J31").Value = 3

# --- Row 32 ---
$ws.Range("C32").Value = 6.5
$ws.Range("J32").Value = 1

# --- Row 34 ---
$ws.Range("C34").Value = 7
$ws.Range("J34").Value = 1

# --- Row 35 ---
$ws.Range("C35").Value = 5
$ws.Range("J35").Value = 3

# --- Column H got wider to fit the new text ---
$ws.Columns("H").ColumnWidth = 86.3

# --- View/selection moved down to show the newly-edited rows ---
$excel.ActiveWindow.ScrollRow = 17
$ws.Range("H37").Select()
